# "Generate Report for Archive"
# - Status text changes from "Ready for handoff" to "In Translation" for the
#   single row tracked on every sheet (Overview!E2:F2, zh-cn!C2, de-de!C2 -
#   they all shared the same string).
# - Because the status text got shorter, the Status column(s) narrow
#   (Overview columns E/F, and column C on the zh-cn / de-de detail sheets).

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E1:F1").ColumnWidth = 12.45

# --- zh-cn detail sheet ----------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C1").ColumnWidth = 12.45

# --- de-de detail sheet ----------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C1").ColumnWidth = 12.45
